$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the newly-tracked "Sun" (column I) hours for each task row
$ws.Range("I8").Value = 4
$ws.Range("I13").Value = 2
$ws.Range("I14").Value = 5
$ws.Range("I15").Value = 2
$ws.Range("I16").Value = 4
$ws.Range("I24").Value = 10
$ws.Range("I25").Value = 8
$ws.Range("I30").Value = 6
$ws.Range("I31").Value = 4
$ws.Range("I33").Value = 3
$ws.Range("I34").Value = 3
$ws.Range("I36").Value = 2
$ws.Range("I37").Value = 3
$ws.Range("I38").Value = 1
$ws.Range("I39").Value = 2
$ws.Range("I40").Value = 6
$ws.Range("I41").Value = 5
$ws.Range("I42").Value = 3
$ws.Range("I43").Value = 2
$ws.Range("I44").Value = 3

# Total row: sum the new column same as the other day columns
$ws.Range("I46").Formula = "=SUM(I5:I44)"

# Move the active selection from H47 to I47
$ws.Range("I47").Select()
